# Apply updated profit/cost figures to each job sheet
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2696.1428
$ws.Range("I111").Value = 1696.5
$ws.Range("K111").Value = 5089.5
$ws.Range("M111").Value = -2022.5
$ws.Range("H113").Value = 4786.825
$ws.Range("I113").Value = 4758.3335
$ws.Range("K113").Value = 4758.3335
$ws.Range("M113").Value = -1504.3335
$ws.Range("H116").Value = 4020.1765
$ws.Range("I116").Value = 3972.6155
$ws.Range("K116").Value = 3972.6155
$ws.Range("M116").Value = -530.6154999999999
$ws.Range("H125").Value = 1330
$ws.Range("I125").Value = 970
$ws.Range("J125").Value = 1450
$ws.Range("K125").Value = 8730
$ws.Range("L125").Value = 13050
$ws.Range("M125").Value = -6270
$ws.Range("N125").Value = -17970
$ws.Range("H132").Value = 799.5833
$ws.Range("I132").Value = 798.13635
$ws.Range("J132").Value = 815.5
$ws.Range("K132").Value = 2394.40905
$ws.Range("L132").Value = 2446.5
$ws.Range("M132").Value = 135.5909499999998
$ws.Range("N132").Value = -7506.5
$ws.Range("H141").Value = 5667.6
$ws.Range("I141").Value = 4715.2856
$ws.Range("J141").Value = 19000
$ws.Range("K141").Value = 14145.8568
$ws.Range("L141").Value = 57000
$ws.Range("M141").Value = -8965.856800000001
$ws.Range("N141").Value = -67360

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2639
$ws.Range("I63").Value = 2639
$ws.Range("K63").Value = 2639
$ws.Range("M63").Value = -1953
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H66").Value = 2639
$ws.Range("I66").Value = 2639
$ws.Range("K66").Value = 13195
$ws.Range("M66").Value = -9763
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H74").Value = 2903.4243
$ws.Range("I74").Value = 2672.2104
$ws.Range("J74").Value = 3217.2144
$ws.Range("K74").Value = 2672.2104
$ws.Range("L74").Value = 3217.2144
$ws.Range("M74").Value = -1798.2104
$ws.Range("N74").Value = -4965.2144
$ws.Range("H77").Value = 2903.4243
$ws.Range("I77").Value = 2672.2104
$ws.Range("J77").Value = 3217.2144
$ws.Range("K77").Value = 13361.052
$ws.Range("L77").Value = 16086.072
$ws.Range("M77").Value = -8993.052
$ws.Range("N77").Value = -24822.072
$ws.Range("H97").Value = 618.86664
$ws.Range("I97").Value = 273.75
$ws.Range("J97").Value = 1999.3334
$ws.Range("K97").Value = 273.75
$ws.Range("L97").Value = 1999.3334
$ws.Range("M97").Value = 222.25
$ws.Range("N97").Value = -2991.3334
$ws.Range("H122").Value = 4049.4783
$ws.Range("I122").Value = 3860.4146
$ws.Range("J122").Value = 5599.8
$ws.Range("K122").Value = 11581.2438
$ws.Range("L122").Value = 16799.4
$ws.Range("M122").Value = -9131.2438
$ws.Range("N122").Value = -21699.4

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7096.25
$ws.Range("J94").Value = 11061
$ws.Range("L94").Value = 11061
$ws.Range("N94").Value = -11963
$ws.Range("H107").Value = 9523.23
$ws.Range("I107").Value = 7840.5
$ws.Range("J107").Value = 10965.571
$ws.Range("K107").Value = 7840.5
$ws.Range("L107").Value = 10965.571
$ws.Range("M107").Value = -5920.5
$ws.Range("N107").Value = -14805.571
$ws.Range("H134").Value = 2668.2546
$ws.Range("I134").Value = 2208.0571
$ws.Range("K134").Value = 6624.1713
$ws.Range("M134").Value = -4089.1713

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3087.8
$ws.Range("I16").Value = 3789.889
$ws.Range("K16").Value = 3789.889
$ws.Range("M16").Value = -3502.889
$ws.Range("H21").Value = 6176
$ws.Range("J21").Value = 6757.5
$ws.Range("L21").Value = 6757.5
$ws.Range("N21").Value = -7227.5
$ws.Range("H29").Value = 2749
$ws.Range("J29").Value = 3498
$ws.Range("L29").Value = 3498
$ws.Range("N29").Value = -4084
$ws.Range("H31").Value = 1467.9
$ws.Range("I31").Value = 1282
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1282
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -987
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 1467.9
$ws.Range("I34").Value = 1282
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 1282
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -1080
$ws.Range("N34").Value = -5404
$ws.Range("H107").Value = 562.94446
$ws.Range("I107").Value = 470
$ws.Range("K107").Value = 470
$ws.Range("M107").Value = 1450
$ws.Range("H113").Value = 3087.8
$ws.Range("I113").Value = 3789.889
$ws.Range("K113").Value = 3789.889
$ws.Range("M113").Value = -1619.889
$ws.Range("H122").Value = 1531
$ws.Range("I122").Value = 1371.8462
$ws.Range("J122").Value = 1760.8889
$ws.Range("K122").Value = 4115.5386
$ws.Range("L122").Value = 5282.6667
$ws.Range("M122").Value = -1665.5386
$ws.Range("N122").Value = -10182.6667
$ws.Range("H132").Value = 1583.4667
$ws.Range("J132").Value = 3994
$ws.Range("L132").Value = 11982
$ws.Range("N132").Value = -17042

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1468.3636
$ws.Range("I5").Value = 1350.3334
$ws.Range("J5").Value = 1999.5
$ws.Range("K5").Value = 4051.0002
$ws.Range("L5").Value = 5998.5
$ws.Range("M5").Value = -3939.0002
$ws.Range("N5").Value = -6222.5
$ws.Range("H12").Value = 1016.7778
$ws.Range("I12").Value = 538.2
$ws.Range("J12").Value = 1615
$ws.Range("K12").Value = 1614.6
$ws.Range("L12").Value = 4845
$ws.Range("M12").Value = -1441.6
$ws.Range("N12").Value = -5191
$ws.Range("H38").Value = 149.44444
$ws.Range("J38").Value = 181.77777
$ws.Range("L38").Value = 545.33331
$ws.Range("N38").Value = -1239.33331
$ws.Range("H59").Value = 9699.875
$ws.Range("I59").Value = 5999.6665
$ws.Range("K59").Value = 17998.9995
$ws.Range("M59").Value = -17458.9995
$ws.Range("H135").Value = 1468.3636
$ws.Range("I135").Value = 1350.3334
$ws.Range("J135").Value = 1999.5
$ws.Range("K135").Value = 12153.0006
$ws.Range("L135").Value = 17995.5
$ws.Range("M135").Value = -9618.000599999999
$ws.Range("N135").Value = -23065.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 4379.4
$ws.Range("I107").Value = 5299.3335
$ws.Range("J107").Value = 2999.5
$ws.Range("K107").Value = 5299.3335
$ws.Range("L107").Value = 2999.5
$ws.Range("M107").Value = -3379.3335
$ws.Range("N107").Value = -6839.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 15000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 15000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 15000
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -15280
$ws.Range("H46").Value = 2175
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2175
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2175
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2551
$ws.Range("H55").Value = 1046.4166
$ws.Range("I55").Value = 585.5
$ws.Range("K55").Value = 585.5
$ws.Range("M55").Value = -412.5
$ws.Range("H61").Value = 3919.1667
$ws.Range("I61").Value = 4000
$ws.Range("J61").Value = 3903
$ws.Range("K61").Value = 4000
$ws.Range("L61").Value = 3903
$ws.Range("M61").Value = -3798
$ws.Range("N61").Value = -4307
$ws.Range("H113").Value = 3919.1667
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 3903
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 3903
$ws.Range("M113").Value = -1830
$ws.Range("N113").Value = -8243

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 9499.5
$ws.Range("I7").Value = 8999
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 8999
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -8886
$ws.Range("N7").Value = -10226
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H20").Value = 43154.57
$ws.Range("I20").Value = 31581.8
$ws.Range("J20").Value = 72086.5
$ws.Range("K20").Value = 31581.8
$ws.Range("L20").Value = 72086.5
$ws.Range("M20").Value = -31341.8
$ws.Range("N20").Value = -72566.5
$ws.Range("H122").Value = 3396.4666
$ws.Range("I122").Value = 3452.861
$ws.Range("K122").Value = 10358.583
$ws.Range("M122").Value = -7908.582999999999
$ws.Range("H132").Value = 1256.9736
$ws.Range("I132").Value = 1132.0938
$ws.Range("J132").Value = 1923
$ws.Range("K132").Value = 3396.2814
$ws.Range("L132").Value = 5769
$ws.Range("M132").Value = -866.2814000000003
$ws.Range("N132").Value = -10829
$ws.Range("H136").Value = 4343.676
$ws.Range("I136").Value = 2944.6177
$ws.Range("K136").Value = 8833.8531
$ws.Range("M136").Value = -6283.8531
